$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Qty"
$ws.Range("B1").Value = "nroPart"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Price"

$ws.Range("A2").Value = 3
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = "Custom Design"
$ws.Range("D2").Value = 75
